$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Cor dos inputs de peso e altura está diferente dos outros" (row 52) as Feito
$ws.Range("D52").Value = "Feito"

# Mark "Input estão com cor incorreta" (row 73) as Feito
$ws.Range("D73").Value = "Feito"

# Remove the resolved item "Select está com estilização diferente" (row 78) entirely,
# shifting subsequent rows up by one
$ws.Rows.Item(78).Delete()

# Mark "Sinalizar campos obrigatórios" (now row 78 after the deletion above) as Feito
$ws.Range("D78").Value = "Feito"

# Update the view to reflect where the user ended up working
$ws.Application.ActiveWindow.ScrollRow = 57
$ws.Range("D80").Select()
